$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "98.311.38"
$ws.Range("E2").Value = "  +4.36%  "

# Row 3
$ws.Range("D3").Value = "3.368.87"
$ws.Range("E3").Value = "  +9.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.12%  "

# Row 7
$ws.Range("E7").Value = "  +8.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.385"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "3.366.11"
$ws.Range("E10").Value = "  +9.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "

# Row 12
$ws.Range("E12").Value = "  +1.21%  "

# Row 13
$ws.Range("D13").Value = "98.009.17"
$ws.Range("E13").Value = "  +4.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.96%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.991.59"
$ws.Range("E15").Value = "  +9.54%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000246"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.40%  "

# Row 17
$ws.Range("E17").Value = "  +3.24%  "

# Row 18
$ws.Range("D18").Value = "3.374.69"
$ws.Range("E18").Value = "  +9.91%  "

# Row 19
$ws.Range("E19").Value = "  +2.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.59%  "

# Row 22
$ws.Range("E22").Value = "  +3.20%  "

# Row 23
$ws.Range("E23").Value = "  +9.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "

# Row 28
$ws.Range("D28").Value = "3.541.32"
$ws.Range("E28").Value = "  +9.39%  "

# Row 29
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.254"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.77%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.187"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.00%  "

# Row 32
$ws.Range("E32").Value = "  +1.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.79%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "527.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.97%  "

# Row 37
$ws.Range("E37").Value = "  +1.40%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.87%  "

# Row 40
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.12%  "

# Row 41
$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.32%  "

# Row 42
$ws.Range("E42").Value = "  +2.93%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.86%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.788"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.33%  "

# Row 46
$ws.Range("E46").Value = "  -0.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.24%  "

# Row 49
$ws.Range("E49").Value = "  +6.79%  "

# Row 50
$ws.Range("E50").Value = "  +4.26%  "

# Row 51
$ws.Range("E51").Value = "  +6.85%  "
